# Insert a new data row at row 70 (pushing existing rows 70-152 down to 71-153)
# and populate it with a new weekly price observation, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 70..152 down by one row, creating a blank row 70 (style of row 70 is
# carried along automatically by Excel's row insert behaviour).
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record.
$ws.Range("A70").Value = 6
$ws.Range("B70").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C70").Value = 'Metropolitana'
$ws.Range("D70").Value = 44494
$ws.Range("E70").Value = 13
$ws.Range("F70").Value = 100112022
$ws.Range("G70").Value = 'Arveja Verde'
$ws.Range("H70").Value = 'Perfection'
$ws.Range("I70").Value = 'Primera'
$ws.Range("J70").Value = 250
$ws.Range("K70").Value = 16000
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = 17200
$ws.Range("N70").Value = '$/malla 25 kilos'
$ws.Range("O70").Value = 'Provincia de Huasco'
$ws.Range("P70").Value = 688
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = 'Hortaliza'
